$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: Insert two new columns before column D.
# This shifts the existing D:K (8 quarters) to F:M, preserving values/styles.
$ws.Columns("D:E").Insert()

# Step 2: Copy number formatting from column F (the shifted former column D)
# onto the two new D:E columns, so date header rows and data rows pick up
# the same styles (s=2 for date headers, s=3 for numeric data) instead of
# the default style Insert() assigns from the left neighbour (column C).
$ws.Range("F7:F102").Copy()
$ws.Range("D7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Step 3: Populate the two new columns (D = Dec-2018 / 43465, E = Sep-2018 / 43373)
# with the latest quarterly figures.
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = "NA"
$ws.Range("E8").Value = 1100
$ws.Range("D9").Value = "NA"
$ws.Range("E9").Value = "NA"
$ws.Range("D10").Value = "NA"
$ws.Range("E10").Value = "NA"
$ws.Range("D12").Value = 13700
$ws.Range("E12").Value = 15600
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 17500
$ws.Range("E17").Value = 19700
$ws.Range("D18").Value = "NA"
$ws.Range("E18").Value = -18600
$ws.Range("D20").Value = "NA"
$ws.Range("E20").Value = 200
$ws.Range("D21").Value = "NA"
$ws.Range("E21").Value = -18100
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("D23").Value = -16900
$ws.Range("E23").Value = -18400
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = -16900
$ws.Range("E26").Value = -18400
$ws.Range("D27").Value = -16900
$ws.Range("E27").Value = -18400
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = "NA"
$ws.Range("E32").Value = -200
$ws.Range("D33").Value = -16900
$ws.Range("E33").Value = -18400
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = -16900
$ws.Range("E35").Value = -18400
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 28800
$ws.Range("E41").Value = 19700
$ws.Range("D42").Value = 89600
$ws.Range("E42").Value = 26700
$ws.Range("D43").Value = "NA"
$ws.Range("E43").Value = "NA"
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("D45").Value = 2500
$ws.Range("E45").Value = 1900
$ws.Range("D46").Value = 120900
$ws.Range("E46").Value = 48300
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 14500
$ws.Range("E48").Value = 14800
$ws.Range("D49").Value = 0
$ws.Range("E49").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 800
$ws.Range("E52").Value = 800
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 136100
$ws.Range("E54").Value = 64000
$ws.Range("D57").Value = 1900
$ws.Range("E57").Value = 2100
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("D59").Value = 6700
$ws.Range("E59").Value = 9400
$ws.Range("D60").Value = 8600
$ws.Range("E60").Value = 11600
$ws.Range("D61").Value = 0
$ws.Range("E61").Value = 0
$ws.Range("D62").Value = 13200
$ws.Range("E62").Value = 13500
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 21800
$ws.Range("E66").Value = 25100
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = -277500
$ws.Range("E72").Value = -260600
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 114300
$ws.Range("E76").Value = 38900
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = -16900
$ws.Range("E81").Value = -18400
$ws.Range("D83").Value = 400
$ws.Range("E83").Value = 300
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = -19300
$ws.Range("E89").Value = -14600
$ws.Range("D91").Value = 0
$ws.Range("E91").Value = 0
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -62600
$ws.Range("E94").Value = 14300
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = 91100
$ws.Range("E100").Value = 100
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = 9100
$ws.Range("E102").Value = -200


# Step 4: A handful of cells in the shifted historical columns were also
# corrected as part of this data refresh (not pure column shifts).
$ws.Range("H21").Value = -13000
$ws.Range("I21").Value = -13600
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("H83").Value = 400
$ws.Range("I83").Value = 300
$ws.Range("F91").Value = "NA"
$ws.Range("G91").Value = "NA"
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = -100
$ws.Range("I94").Value = 8700
$ws.Range("I102").Value = -4100

"done"
